# Auto-generated edit script: update cryptos price/volume table
# Mirrors the commit "Updated cryptos list on Wed Apr  3 17:50:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '66.003.06'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '3.325.56'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '187.96'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('D6').Value = '553.93'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.316.30'
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').Value = '0.579'
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('E10').Value = '  -3.61%  '
$ws.Range('D11').Value = '0.580'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').Value = '45.93'
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('D13').Value = '0.0000266'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').Value = '3.851.50'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').Value = '8.47'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').Value = '579.82'
$ws.Range('E16').Value = '  -7.95%  '
$ws.Range('D17').Value = '66.004.09'
$ws.Range('E17').Value = '  +0.96%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.117'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.324.99'
$ws.Range('E19').Value = '  +2.16%  '
$ws.Range('D20').Value = '17.78'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('D21').Value = '10.90'
$ws.Range('E21').Value = '  -3.53%  '
$ws.Range('D22').Value = '0.893'
$ws.Range('E22').Value = '  -0.64%  '
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('D24').Value = '5.01'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('D25').Value = '99.11'
$ws.Range('E25').Value = '  -6.12%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').Value = '2.69'
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').Value = '9.33'
$ws.Range('E28').Value = '  -1.77%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '30.55'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '8.40'
$ws.Range('E30').Value = '  -2.80%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '6.59'
$ws.Range('E31').Value = '  +5.21%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '575.55'
$ws.Range('E32').Value = '  +5.48%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D33').Value = '3.76'
$ws.Range('E33').Value = '  -5.49%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = '10.86'
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.103'
$ws.Range('E35').Value = '  -1.06%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.720.97'
$ws.Range('E36').Value = '  +3.16%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').Value = '55.49'
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '34.08'
$ws.Range('E39').Value = '  +6.87%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.126'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0690'
$ws.Range('E41').Value = '  -3.16%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.63'
$ws.Range('E42').Value = '  -3.42%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '3.14'
$ws.Range('E43').Value = '  -7.00%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').Value = '3.36'
$ws.Range('E44').Value = '  +2.49%  '
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').Value = '0.335'
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0409'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('B47').Value = 'CoreDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D47').Value = '2.97'
$ws.Range('E47').Value = '  -8.69%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('B49').Value = 'FirstDigitalUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D49').Value = '0.999'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '2.52'
$ws.Range('E50').Value = '  -3.04%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '126.62'
$ws.Range('E51').Value = '  +5.92%  '
